# Applies the "effects" benchmark results update:
#  - adds a second results table (rows 21-28) below the existing one
#  - adds a new empty, italic/gray "placeholder" cell at B33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Second table header (row 21) - reuse the bold/centered style already
#    used by the first table's header row (row 2) so no new styles/fonts
#    are introduced.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B21:E21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B21").Value = "Program"
$ws.Range("C21").Value = "LoC"
$ws.Range("D21").Value = "LoS"
$ws.Range("E21").Value = "Time(ms)"

# ---------------------------------------------------------------------------
# 2. Second table body (rows 22-28) - reuse the plain centered style
#    already used in column F (e.g. F3) for the new cells.
# ---------------------------------------------------------------------------
$ws.Range("F3").Copy() | Out-Null
$ws.Range("B22:E28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$rows = @(
    @{ Row = 22; Id = 1; Name = "State Monad "     ; LoC = 29 ; LoS = 6   },
    @{ Row = 23; Id = 2; Name = "Exchange State"    ; LoC = 19 ; LoS = 3   },
    @{ Row = 24; Id = 3; Name = "Inductive sumEff"  ; LoC = 20 ; LoS = 5   },
    @{ Row = 25; Id = 4; Name = "Generic_Count"     ; LoC = 32 ; LoS = 7   },
    @{ Row = 26; Id = 5; Name = "McCarthy_Andrej_Bauer"; LoC = 83 ; LoS = $null },
    @{ Row = 27; Id = 6; Name = "McCarthy_Pythagorean" ; LoC = 70 ; LoS = $null },
    @{ Row = 28; Id = 7; Name = "Async Yield"       ; LoC = 109; LoS = $null }
)

# Write the "Name" column first, in the same order the new labels were
# originally authored (State Monad, Inductive sumEff, Async Yield,
# Exchange State), so brand-new shared-string entries land in that order.
$nameOrder = @(22, 24, 28, 23, 25, 26, 27)
foreach ($rowNum in $nameOrder) {
    $r = $rows | Where-Object { $_.Row -eq $rowNum }
    $ws.Cells.Item($r.Row, 2).Value = $r.Name
}

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Id
    $ws.Cells.Item($r.Row, 3).Value = $r.LoC
    if ($null -ne $r.LoS) {
        $ws.Cells.Item($r.Row, 4).Value = $r.LoS
    }
}

# ---------------------------------------------------------------------------
# 3. New trailing placeholder cell (B33), italic light-gray "Menlo" font,
#    left blank - used as an editing cursor position marker.
#
#    Building the font via a transient named style (removed again right
#    after use) keeps the generated font table clean - it avoids leaving
#    behind the extra intermediate font variations that accumulate when
#    mutating a Range's Font object property-by-property.
# ---------------------------------------------------------------------------
$tempStyleName = "__tmp_placeholder_style__"
$tempStyle = $wb.Styles.Add($tempStyleName)
$tempStyle.Font.Name = "Menlo"
$tempStyle.Font.Size = 12
$tempStyle.Font.Italic = $true
$tempStyle.Font.Color = 11184810   # RGB(170,170,170) -> FFAAAAAA

$placeholder = $ws.Range("B33")
$placeholder.Style = $tempStyleName
$wb.Styles.Item($tempStyleName).Delete()

# Move the active selection to the new placeholder cell, matching the
# author's final cursor position.
$placeholder.Select() | Out-Null
